# Merging of suites and updation of code
# Update the OrderDate / OverageID test-data row on "Sheet1" to the latest
# values used by the automation suite.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A ("OrderDate") -> 10-28-2021
$ws.Range("A2").Formula = '="10-28-2021"'
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

# Column L ("OverageID") -> 57905035
$ws.Range("L2").Formula = '="57905035"'
$ws.Range("L2").Copy()
$ws.Range("L2").PasteSpecial(-4163)

$excel.CutCopyMode = $false
